$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for "CCP" (Concepcion, Chile) - row 177.
# All rows below it (178:302) shift up by one.
$ws.Rows.Item(177).Delete()
